$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "scenario" column header in F1, matching the bold/centered
# header style used by the other header cells (A1:E1).
$ws.Range("F1").Value = "scenario"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108

# Fill F2:F101 with the scenario label "S3" for every data row.
$lastRow = 101
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = "S3"
}
